$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.201.58'
$ws.Range('E2').Value = '  -3.84%  '
$ws.Range('D3').Value = '1.964.98'
$ws.Range('E3').Value = '  -6.51%  '
$ws.Range('E4').Value = '  +1.76%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.02'
$ws.Range('E5').Value = '  -4.53%  '
$ws.Range('E6').Value = '  +1.84%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4991'
$ws.Range('E7').Value = '  -6.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4196'
$ws.Range('E8').Value = '  -5.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.45'
$ws.Range('E9').Value = '  -2.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08865'
$ws.Range('E10').Value = '  -5.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.097'
$ws.Range('E11').Value = '  -6.39%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.97'
$ws.Range('E12').Value = '  -7.41%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '2.057.37'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.871'
$ws.Range('E14').Value = '  -8.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.411'
$ws.Range('E15').Value = '  -7.43%  '
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001100'
$ws.Range('E17').Value = '  -5.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.50'
$ws.Range('E18').Value = '  -10.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06709'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.22'
$ws.Range('E20').Value = '  -9.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.017'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.988'
$ws.Range('E22').Value = '  -5.52%  '
$ws.Range('D23').Value = '29.333.71'
$ws.Range('E23').Value = '  -3.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.89'
$ws.Range('E24').Value = '  -5.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.313'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.60'
$ws.Range('E26').Value = '  -5.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.92'
$ws.Range('E27').Value = '  -4.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.191'
$ws.Range('E28').Value = '  -8.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.285'
$ws.Range('E29').Value = '  -9.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.61'
$ws.Range('E30').Value = '  -5.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.046'
$ws.Range('E31').Value = '  -8.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09858'
$ws.Range('E32').Value = '  -6.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.499'
$ws.Range('E33').Value = '  -10.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.783'
$ws.Range('E34').Value = '  -7.65%  '
$ws.Range('E35').Value = '  -3.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02431'
$ws.Range('E36').Value = '  -8.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.224'
$ws.Range('E37').Value = '  -9.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.289'
$ws.Range('E38').Value = '  -4.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06315'
$ws.Range('E39').Value = '  -7.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6457'
$ws.Range('E40').Value = '  -8.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.52'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2016'
$ws.Range('E42').Value = '  -9.47%  '
$ws.Range('E43').Value = '  +1.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6235'
$ws.Range('E44').Value = '  -9.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.49'
$ws.Range('E45').Value = '  -6.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.185'
$ws.Range('E46').Value = '  -6.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.285'
$ws.Range('E47').Value = '  -7.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.480'
$ws.Range('E48').Value = '  -4.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000335'
$ws.Range('E49').Value = '  -4.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06901'
$ws.Range('E50').Value = '  -4.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.113'
$ws.Range('E51').Value = '  -9.22%  '
